# Correct values for biogenic methane.
#
# The "Methane, non-fossil" block (rows 192-200) had 9 rows with
# duplicated / inconsistent category combinations. It should instead
# have exactly 5 rows (one per category, same pattern used for every
# other substance), with amounts computed as the corresponding
# "Methane, fossil" amount (rows 177-181) minus 2.75.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the 4 extra / duplicate rows (the old rows 197-200) so the
# "Methane, non-fossil" block shrinks from 9 rows down to 5 (192-196),
# matching the 5-category pattern used by every other substance.
$ws.Range("A197:C200").EntireRow.Delete()

# Re-derive each of the 5 remaining "Methane, non-fossil" amounts from
# the matching "Methane, fossil" row (same category, 15 rows above).
$ws.Cells.Item(192, 3).Formula = "=C177-2.75"
$ws.Cells.Item(193, 3).Formula = "=C178-2.75"
$ws.Cells.Item(194, 3).Formula = "=C179-2.75"
$ws.Cells.Item(195, 3).Formula = "=C180-2.75"
$ws.Cells.Item(196, 3).Formula = "=C181-2.75"

# The data range shrank from A1:C253 to A1:C249 - refresh the
# AutoFilter and the hidden _FilterDatabase name to match.
if ($ws.AutoFilterMode) {
    $ws.AutoFilterMode = $false
}
$ws.Range("A1:C249").AutoFilter()

foreach ($n in $wb.Names) {
    if ($n.Name() -eq "Sheet1!_FilterDatabase") {
        $n.RefersTo = "=Sheet1!`$A`$1:`$C`$249"
    }
}

$wb.Application.Calculate()
